$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("Include ValueSets").Name = "Include ValueSet #0"
$wb.Worksheets.Item("Include from CareSocialCodes").Name = "Include #1"

# Update Metadata values
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("B3").Value = "2.1.0"
$ws.Range("B8").Value = "2024-10-31T20:37:15+01:00"
